# Auto-generated edit script applying the bilibili-scraped convention update
# (commit: 'Update gh-pages to output generated at 456a3b4')
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # 展览
$ws2 = $wb.Worksheets.Item(2)  # 演出
$ws3 = $wb.Worksheets.Item(3)  # 本地生活
$ws4 = $wb.Worksheets.Item(4)  # 全部类型


# ---- 展览 (sheet1) ----
$ws1.Range("F2").Value = 1939
$ws1.Range("F5").Value = 417
$ws1.Range("F6").Value = 1831
$ws1.Range("F8").Value = 1299
$ws1.Range("F9").Value = 522
$ws1.Range("F11").Value = 2732
$ws1.Range("F12").Value = 366
$ws1.Range("F14").Value = 1079
$ws1.Range("F15").Value = 577
$ws1.Range("F16").Value = 22
$ws1.Range("F17").Value = 58
$ws1.Range("F18").Value = 1572
$ws1.Range("F20").Value = 1234
$ws1.Range("F21").Value = 174
$ws1.Range("C23").Value = "上海·坏孩纸物语の第48届动漫节之梦回春秋战国（免费活动）"
$ws1.Range("D23").Value = "世纪大道2002号 S.C.Plaza"
$ws1.Range("E23").Value = "2024.07.13 10:00-07.14 17:00"
$ws1.Range("F23").Value = 2
$ws1.Range("G23").Value = 49.6
$ws1.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=88004"
$ws1.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202406/PP9QKg0v1719293500048.png"
$ws1.Range("F24").Value = 52
$ws1.Range("F25").Value = 1439
$ws1.Range("F26").Value = 1448
$ws1.Range("F27").Value = 1316
$ws1.Range("F28").Value = 235
$ws1.Range("F29").Value = 1270
$ws1.Range("F31").Value = 143
$ws1.Range("F34").Value = 1831
$ws1.Range("F35").Value = 464
$ws1.Range("F37").Value = 150
$ws1.Range("F39").Value = 2256
$ws1.Range("F40").Value = 137
$ws1.Range("F42").Value = 2747
$ws1.Range("F45").Value = 16

# ---- 演出 (sheet2) ----
$ws2.Range("F5").Value = 58
$ws2.Range("F12").Value = 360
$ws2.Range("F13").Value = 108111
$ws2.Range("F17").Value = 62
$ws2.Range("F18").Value = 62
$ws2.Range("F22").Value = 274
$ws2.Range("F23").Value = 66
$ws2.Range("F30").Value = 31
$ws2.Range("F37").Value = 161

# ---- 本地生活 (sheet3) ----
$ws3.Range("F5").Value = 3016
$ws3.Range("F6").Value = 4846
$ws3.Range("G6").Value = "已售罄"
$ws3.Range("F10").Value = 926
$ws3.Range("F12").Value = 600
$ws3.Range("F13").Value = 1339
$ws3.Range("F14").Value = 377
$ws3.Range("F15").Value = 1195

# ---- 全部类型 (sheet4) ----
$ws4.Range("F2").Value = 1939
$ws4.Range("B5").NumberFormat = "@"
$ws4.Range("B5").Value = "2024-06-01"
$ws4.Range("B5").Style = "Normal"
$ws4.Range("C5").Value = "上海·NIJISANJI EN 官方授权主题店"
$ws4.Range("D5").Value = "西藏北路166号（地铁8号线曲阜路下） 静安大悦城"
$ws4.Range("E5").Value = "2024.06.01 00:00-07.15 23:59"
$ws4.Range("F5").Value = 655
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=86310"
$ws4.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202405/MhBVkfZ51716778059321.jpeg"
$ws4.Range("B6").NumberFormat = "@"
$ws4.Range("B6").Value = "2024-06-07"
$ws4.Range("B6").Style = "Normal"
$ws4.Range("C6").Value = "上海·全职高手×HAPPY ZOO 全职高手十周年咖啡厅"
$ws4.Range("D6").Value = "南京东路340号百联zx创趣场四楼05号 HAPPY ZOO"
$ws4.Range("E6").Value = "2024.06.07 00:00-08.04 23:59"
$ws4.Range("F6").Value = 926
$ws4.Range("G6").Value = 10
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=86871"
$ws4.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202406/KLJmCEkC1717568198482.png"
$ws4.Range("B7").NumberFormat = "@"
$ws4.Range("B7").Value = "2024-06-08"
$ws4.Range("B7").Style = "Normal"
$ws4.Range("C7").Value = "上海· 怪兽8号 meets niko and … 集章之旅    "
$ws4.Range("D7").Value = "吴江路169号1层E127,E128 niko and ... (上海四季坊店)"
$ws4.Range("E7").Value = "2024.06.08 10:00-07.21 22:00"
$ws4.Range("F7").Value = 534
$ws4.Range("G7").Value = 48
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=85758"
$ws4.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202405/xw8aUE5u1715846379865.jpeg"
$ws4.Range("C8").Value = "上海·无穹-中国 航天沉浸艺术展"
$ws4.Range("D8").Value = "上海浦东新区樱花路869号3F 上海喜玛拉雅美术馆"
$ws4.Range("E8").Value = "2024.06.08 10:00-10.07 20:00"
$ws4.Range("F8").Value = 96
$ws4.Range("G8").Value = 78
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=86957"
$ws4.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202406/Bus3lAnI1717558639134.jpeg"
$ws4.Range("B9").NumberFormat = "@"
$ws4.Range("B9").Value = "2024-06-12"
$ws4.Range("B9").Style = "Normal"
$ws4.Range("C9").Value = "上海 ·「蓝色监狱 x 次元波板糖」主题餐厅"
$ws4.Range("D9").Value = "西藏北路166号（地铁8号线曲阜路下） 静安大悦城"
$ws4.Range("E9").Value = "2024.06.12 00:00-06.30 23:59"
$ws4.Range("F9").Value = 600
$ws4.Range("G9").Value = 30
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=87144"
$ws4.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202406/budG5Uyl1718089258239.png"
$ws4.Range("B10").NumberFormat = "@"
$ws4.Range("B10").Value = "2024-06-14"
$ws4.Range("B10").Style = "Normal"
$ws4.Range("C10").Value = "上海·「排球少年!!垃圾场决战」主题店"
$ws4.Range("E10").Value = "2024.06.14 00:00-07.07 23:59"
$ws4.Range("F10").Value = 1339
$ws4.Range("G10").Value = 10
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=86948"
$ws4.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202406/GxQRrJ2V1717655942245.png"
$ws4.Range("B11").NumberFormat = "@"
$ws4.Range("B11").Value = "2024-06-29"
$ws4.Range("B11").Style = "Normal"
$ws4.Range("C11").Value = "上海·cdc动漫展"
$ws4.Range("D11").Value = "海潮路133号B1 JUMP工坊"
$ws4.Range("E11").Value = "2024.06.29 10:00-06.30 17:00"
$ws4.Range("F11").Value = 418
$ws4.Range("G11").Value = 60
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=85110"
$ws4.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202405/RMpaP6sF1714725969882.jpeg"
$ws4.Range("C12").Value = "上海·创世次元动漫游戏嘉年华3.0"
$ws4.Range("D12").Value = "中环立交桥苏宁天御国际广场西南侧约240米 轮客行轮滑馆(普陀店)"
$ws4.Range("F12").Value = 1832
$ws4.Range("G12").Value = 58
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=86506"
$ws4.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202405/Clkfdwic1716894666596.jpeg"
$ws4.Range("C13").Value = "上海·第六十三届燃梦星辰动漫嘉年华"
$ws4.Range("D13").Value = "陆宝山路155号 佘山·旭辉里"
$ws4.Range("E13").Value = "2024.06.29 14:00-06.29 18:00"
$ws4.Range("F13").Value = 857
$ws4.Range("G13").Value = 58.8
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=85231"
$ws4.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202406/Tx1g80PC1717667546414.jpeg"
$ws4.Range("F14").Value = 1299
$ws4.Range("F16").Value = 522
$ws4.Range("F17").Value = 1195
$ws4.Range("F18").Value = 2732
$ws4.Range("F20").Value = 366
$ws4.Range("F22").Value = 1079
$ws4.Range("F23").Value = 577
$ws4.Range("F24").Value = 22
$ws4.Range("F25").Value = 1572
$ws4.Range("F27").Value = 360
$ws4.Range("F28").Value = 1234
$ws4.Range("F29").Value = 174
$ws4.Range("F31").Value = 1439
$ws4.Range("F32").Value = 1448
$ws4.Range("F33").Value = 1316
$ws4.Range("F35").Value = 62
$ws4.Range("F36").Value = 1270
$ws4.Range("F39").Value = 66
$ws4.Range("F40").Value = 1831
$ws4.Range("F43").Value = 2256
$ws4.Range("F44").Value = 137
$ws4.Range("F46").Value = 2747
